$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date strings (slashes -> dashes) for each row.
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    # Some of the new day-first strings (day <= 12) are ambiguous and would
    # otherwise get auto-parsed into a real date serial by Excel's smart
    # entry. Force the cell to Text first so the literal string is kept,
    # then restore the "Normal" style so no stray number-format sticks to
    # the cell (matches the original plain-text cell formatting).
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.Style = "Normal"
}

# Row 3 (28-07-2022): Total Attendance Count (D) and Invalid (G) flip 0 -> 1.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Rows 7, 13, 16, 20: Total Attendance Count (D) and Real (E) flip 0 -> 1,
# Absent (H) flips 1 -> 0.
$rowsDEH = @(7, 13, 16, 20)
foreach ($row in $rowsDEH) {
    $ws.Range("D$row").Value = 1
    $ws.Range("E$row").Value = 1
    $ws.Range("H$row").Value = 0
}
